$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N ("Late"), pushing the
# "Late" / "Original"(heading) / "Outstanding" columns one place to the
# right (N->O, O->P, P->Q). The new column inherits the width of the
# column immediately to its left (M), matching Excel's normal insert
# behaviour.
$ws.Columns("N").Insert() | Out-Null
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Make "Repayment schedule" the active sheet/tab, and select cell K18 on
# it (this also clears tabSelected on whichever sheet had it before).
$ws.Activate() | Out-Null
$ws.Range("K18").Select() | Out-Null

Write-Output "done"
